# Applies the "normative rule linking / adoc formatting" update to the
# test-norm-rules workbook:
#   * Reworks the "no_tag" rule description (adds a link + adoc formatting
#     examples).
#   * Adds "[.underline]#description#" formatting to the
#     paragraph-with-a-really-wide-rule-name description.
#   * Inserts two new rule rows ("bold" and "italics") right after
#     "double_tags" (before "superscript").
#   * Updates the superscript/subscript rule descriptions.
#   * Grows the Table1 ListObject (and therefore the sheet dimension) to
#     account for the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "no_tag" rule description (row 3, column C) ------------
$ws.Range("C3").Value2 = 'Normative rule *without* tag/tags
This normative rule has no references to the standard. This should only be used in extraordinary circumstances.
It does include a link to <<table1>> (another normative rule).
Has basic adoc formatting such as *bold*, ita__lics__, `monospace`, 2^superscript^, ~subscript~, [.underline]#underline#,
and &le; (Unicode text for less-than-equals-to) and &#8800; (Unicode decimal value for not-equal-to).'

# --- 2. Update the "paragraph-with-a-really-wide-rule-name" description ---
$ws.Range("C5").Value2 = 'Here''s a [.underline]#description#.
It''s got 2 lines.
Paragraph without inline anchors'

# --- 3. Insert two new rows for the "bold" and "italics" rules ------------
# They land right after "double_tags" (row 10) and before "superscript"
# (row 11), pushing superscript/subscript and everything below down by two
# rows.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(12).Insert()

$ws.Range("A11").Value2 = "my-chapter_name"
$ws.Range("B11").Value2 = "bold"
$ws.Range("C11").Value2 = "ABC is a network - Bold is removed by tags backend so I don't see it"
$ws.Range("D11").Value2 = '["norm:bold"]'

$ws.Range("A12").Value2 = "my-chapter_name"
$ws.Range("B12").Value2 = "italics"
$ws.Range("C12").Value2 = "Let's have fun today - Italics is removed by tags backend so I don't see it"
$ws.Range("D12").Value2 = '["norm:italics"]'

# --- 4. Update the superscript/subscript descriptions (now rows 13/14) ----
$ws.Range("C13").Value2 = "both 2^32^ and ^32^ work"
$ws.Range("C14").Value2 = "both ~log~ and log~2~ work"

# --- 5. Grow the table to cover the two newly-inserted rows ---------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F51"))
